$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 102: keep Emisora (A) as-is, change Serie (B) and Descripción (C)
# to match the new, smaller "Titulos Referenciados..." catalogue block.
$ws.Range("A102").Value = "SDHMX"
$ws.Range("B102").Value = "N"
$ws.Range("C102").Value = "Titulos Referenciados a acciones Tracs Extranjeras (Tipo de cambio Spot)"

# Update row 103 the same way, with its own Emisora.
$ws.Range("A103").Value = "VMSTX"
$ws.Range("B103").Value = "N"
$ws.Range("C103").Value = "Titulos Referenciados a acciones Tracs Extranjeras (Tipo de cambio Spot)"

# Remove the old "Acciones de Sociedades de Inversion" rows (104-117) that
# are no longer part of the catalogue.
$ws.Range("A104:E117").EntireRow.Delete()

# Restore the selection to what the author left it at on save.
$ws.Range("C71").Select()
